# 1.5.1.xlsx — add a 2021 data column (column R) to the table on sheet1,
# mirroring the formatting of the existing 2020 column (Q) row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R values, keyed by row number. "-" cells reuse the same
# shared "-" placeholder string already used elsewhere in column Q.
$data = [ordered]@{
    4  = 2021
    5  = 109
    6  = 74
    7  = 35
    8  = 36
    9  = 35
    10 = 1
    11 = 15
    12 = 8
    13 = 7
    14 = 12
    15 = 7
    16 = 5
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 17
    21 = 8
    22 = 9
    23 = 9
    24 = 7
    25 = 2
    26 = 20
    27 = 9
    28 = 11
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

# Row 3 is the blank, bottom-bordered spacer row above the header -
# give R3 the same border-only formatting as Q3, no value.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)

foreach ($row in $data.Keys) {
    $ws.Range("Q$row").Copy()
    $ws.Range("R$row").PasteSpecial(-4122)
    $ws.Range("R$row").Value = $data[$row]
}

$excel.CutCopyMode = $false

# Reset the active selection to A1 (the saved workbook no longer parks the
# selection on R13 once the new column is populated).
$ws.Range("A1").Select()
